{"js": "// Apply the text revisions described by the commit:\n// title/meta copy rewritten, and several \"what we like\" / \"what we\n// don't like\" bullet points reworded.\nconst body = context.document.body;\n\nconst replacements = [\n  [\n    \"Play Family Guy for Free - Slot Game Review\",\n    \"Play Family Guy Slot Free\",\n  ],\n  [\n    \"Fantastic graphics and design\",\n    \"Top-of-the-line graphics\",\n  ],\n  [\n    \"Innovative special functions\",\n    \"Special functions add excitement\",\n  ],\n  [\n    \"Opportunities for high payouts\",\n    \"Opportunity for high winnings\",\n  ],\n  [\n    \"May be difficult to understand at first\",\n    \"Minimalist design may be difficult to understand at first\",\n  ],\n  [\n    \"Minimalist design could be improved\",\n    \"Limited similar game options\",\n  ],\n  [\n    \"Discover the fantastic graphics, innovative functions and opportunities for high payouts in our review of Family Guy slot game. Play for free now.\",\n    \"Read our review of Family Guy slot machine and play for free. Exciting gameplay and high winnings.\",\n  ],\n];\n\nfor (const [find, replace] of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the text revisions described by the commit:\n# title/meta copy rewritten, and several \"what we like\" / \"what we\n# don't like\" bullet points reworded.\n$d = $word.ActiveDocument\n\n# Longer / more specific strings are replaced before the shorter bullet\n# fragments they contain (case-sensitive matching below also guards\n# against the lowercase \"opportunities for high payouts\" inside the meta\n# description colliding with the \"Opportunities for high payouts\" bullet).\n$replacements = @(\n    @(\"Discover the fantastic graphics, innovative functions and opportunities for high payouts in our review of Family Guy slot game. Play for free now.\", \"Read our review of Family Guy slot machine and play for free. Exciting gameplay and high winnings.\"),\n    @(\"Play Family Guy for Free - Slot Game Review\", \"Play Family Guy Slot Free\"),\n    @(\"Fantastic graphics and design\", \"Top-of-the-line graphics\"),\n    @(\"Innovative special functions\", \"Special functions add excitement\"),\n    @(\"Opportunities for high payouts\", \"Opportunity for high winnings\"),\n    @(\"May be difficult to understand at first\", \"Minimalist design may be difficult to understand at first\"),\n    @(\"Minimalist design could be improved\", \"Limited similar game options\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.MatchCase = $true\n    $find.Replacement.Text = $replaceText\n    $find.Execute($findText, $true, $false, $false, $false, $false, $true, \"wdFindContinue\", $false, $replaceText, \"wdReplaceAll\")\n}\n"}
